$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 361
$ws.Range("I9").Value = 387.3
$ws.Range("J9").Value = 273.33334
$ws.Range("K9").Value = 387.3
$ws.Range("L9").Value = 273.33334
$ws.Range("M9").Value = -218.3
$ws.Range("N9").Value = -611.33334
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H33").Value = 382.78262
$ws.Range("I33").Value = 215.36842
$ws.Range("K33").Value = 215.36842
$ws.Range("M33").Value = 13.63158000000001
$ws.Range("H58").Value = 692.4167
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9300
$ws.Range("H61").Value = 929.5
$ws.Range("I61").Value = 929.5
$ws.Range("K61").Value = 2788.5
$ws.Range("M61").Value = -2616.5
$ws.Range("H82").Value = 4874.75
$ws.Range("I82").Value = 4874.75
$ws.Range("K82").Value = 14624.25
$ws.Range("M82").Value = -14218.25
$ws.Range("H85").Value = 4874.75
$ws.Range("I85").Value = 4874.75
$ws.Range("K85").Value = 14624.25
$ws.Range("M85").Value = -13220.25
$ws.Range("H87").Value = 97991.8
$ws.Range("J87").Value = 99989.75
$ws.Range("L87").Value = 99989.75
$ws.Range("N87").Value = -102485.75
$ws.Range("H90").Value = 97991.8
$ws.Range("J90").Value = 99989.75
$ws.Range("L90").Value = 299969.25
$ws.Range("N90").Value = -312449.25
$ws.Range("H101").Value = 1564.9
$ws.Range("I101").Value = 1181.125
$ws.Range("K101").Value = 3543.375
$ws.Range("M101").Value = -1921.375
$ws.Range("H132").Value = 56918.215
$ws.Range("I132").Value = 58292.953
$ws.Range("K132").Value = 174878.859
$ws.Range("M132").Value = -172348.859
$ws.Range("H141").Value = 2929.5557
$ws.Range("I141").Value = 2929.5557
$ws.Range("K141").Value = 8788.667099999999
$ws.Range("M141").Value = -3608.667099999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1843.5625
$ws.Range("I45").Value = 1941.5
$ws.Range("J45").Value = 1549.75
$ws.Range("K45").Value = 1941.5
$ws.Range("L45").Value = 1549.75
$ws.Range("M45").Value = -1564.5
$ws.Range("N45").Value = -2303.75
$ws.Range("H57").Value = 11121277
$ws.Range("I57").Value = 11121277
$ws.Range("K57").Value = 11121277
$ws.Range("M57").Value = -11120793
$ws.Range("H63").Value = 10235.714
$ws.Range("I63").Value = 1829.3
$ws.Range("K63").Value = 1829.3
$ws.Range("M63").Value = -1143.3
$ws.Range("H66").Value = 10235.714
$ws.Range("I66").Value = 1829.3
$ws.Range("K66").Value = 9146.5
$ws.Range("M66").Value = -5714.5
$ws.Range("H74").Value = 10273.25
$ws.Range("I74").Value = 10037.467
$ws.Range("J74").Value = 10666.223
$ws.Range("K74").Value = 10037.467
$ws.Range("L74").Value = 10666.223
$ws.Range("M74").Value = -9163.467000000001
$ws.Range("N74").Value = -12414.223
$ws.Range("H77").Value = 10273.25
$ws.Range("I77").Value = 10037.467
$ws.Range("J77").Value = 10666.223
$ws.Range("K77").Value = 50187.33500000001
$ws.Range("L77").Value = 53331.115
$ws.Range("M77").Value = -45819.33500000001
$ws.Range("N77").Value = -62067.115
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4538.8
$ws.Range("I99").Value = 4841.2856
$ws.Range("K99").Value = 4841.2856
$ws.Range("M99").Value = -3343.2856
$ws.Range("H128").Value = 5030
$ws.Range("I128").Value = 5030
$ws.Range("K128").Value = 15090
$ws.Range("M128").Value = -12600
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 76458.7
$ws.Range("I31").Value = 137021.34
$ws.Range("J31").Value = 18419.5
$ws.Range("K31").Value = 137021.34
$ws.Range("L31").Value = 18419.5
$ws.Range("M31").Value = -136726.34
$ws.Range("N31").Value = -19009.5
$ws.Range("H34").Value = 76458.7
$ws.Range("I34").Value = 137021.34
$ws.Range("J34").Value = 18419.5
$ws.Range("K34").Value = 137021.34
$ws.Range("L34").Value = 18419.5
$ws.Range("M34").Value = -136819.34
$ws.Range("N34").Value = -18823.5
$ws.Range("H76").Value = 4869
$ws.Range("I76").Value = 4869
$ws.Range("K76").Value = 4869
$ws.Range("M76").Value = -4554
$ws.Range("H79").Value = 4869
$ws.Range("I79").Value = 4869
$ws.Range("K79").Value = 4869
$ws.Range("M79").Value = -3777
$ws.Range("H122").Value = 4556.125
$ws.Range("I122").Value = 3749.75
$ws.Range("J122").Value = 5362.5
$ws.Range("K122").Value = 11249.25
$ws.Range("L122").Value = 16087.5
$ws.Range("M122").Value = -8799.25
$ws.Range("N122").Value = -20987.5
$ws.Range("H132").Value = 33343920
$ws.Range("I132").Value = 52639144
$ws.Range("J132").Value = 15803.909
$ws.Range("K132").Value = 157917432
$ws.Range("L132").Value = 47411.727
$ws.Range("M132").Value = -157914902
$ws.Range("N132").Value = -52471.727
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 235.10527
$ws.Range("J2").Value = 71.333336
$ws.Range("L2").Value = 428.000016
$ws.Range("N2").Value = -654.000016
$ws.Range("H17").Value = 126.14286
$ws.Range("J17").Value = 625
$ws.Range("L17").Value = 1875
$ws.Range("N17").Value = -2213
$ws.Range("H33").Value = 249.5
$ws.Range("I33").Value = 115.42857
$ws.Range("J33").Value = 437.2
$ws.Range("K33").Value = 692.57142
$ws.Range("L33").Value = 2623.2
$ws.Range("M33").Value = -409.57142
$ws.Range("N33").Value = -3189.2
$ws.Range("H38").Value = 61.285713
$ws.Range("I38").Value = 64.833336
$ws.Range("K38").Value = 194.500008
$ws.Range("M38").Value = 152.499992
$ws.Range("H103").Value = 941
$ws.Range("I103").Value = 1200.8334
$ws.Range("J103").Value = 421.33334
$ws.Range("K103").Value = 3602.5002
$ws.Range("L103").Value = 1264.00002
$ws.Range("M103").Value = -2723.5002
$ws.Range("N103").Value = -3022.00002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 5000
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 5000
$ws.Range("N96").Value = -10492
$ws.Range("H117").Value = 48128
$ws.Range("J117").Value = 48128
$ws.Range("L117").Value = 48128
$ws.Range("N117").Value = -55012
$ws.Range("H132").Value = 21541016
$ws.Range("I132").Value = 30673836
$ws.Range("J132").Value = 13654.5
$ws.Range("K132").Value = 92021508
$ws.Range("L132").Value = 40963.5
$ws.Range("M132").Value = -92018978
$ws.Range("N132").Value = -46023.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1578
$ws.Range("I46").Value = 964.5714
$ws.Range("K46").Value = 964.5714
$ws.Range("M46").Value = -776.5714
$ws.Range("H120").Value = 92358.60000000001
$ws.Range("J120").Value = 92358.60000000001
$ws.Range("L120").Value = 92358.60000000001
$ws.Range("N120").Value = -102034.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H70").Value = 50105
$ws.Range("J70").Value = 50105
$ws.Range("L70").Value = 50105
$ws.Range("N70").Value = -50735
$ws.Range("H73").Value = 50105
$ws.Range("J73").Value = 50105
$ws.Range("L73").Value = 50105
$ws.Range("N73").Value = -52289
$ws.Range("H74").Value = 197357.17
$ws.Range("I74").Value = 40485.75
$ws.Range("K74").Value = 40485.75
$ws.Range("M74").Value = -39549.75
$ws.Range("H77").Value = 197357.17
$ws.Range("I77").Value = 40485.75
$ws.Range("K77").Value = 121457.25
$ws.Range("M77").Value = -116777.25
$ws.Range("H81").Value = 9624.75
$ws.Range("I81").Value = 7499.5
$ws.Range("K81").Value = 14999
$ws.Range("M81").Value = -13938
$ws.Range("H84").Value = 9624.75
$ws.Range("I84").Value = 7499.5
$ws.Range("K84").Value = 74995
$ws.Range("M84").Value = -69691
$ws.Range("H123").Value = 74995
$ws.Range("J123").Value = 74995
$ws.Range("L123").Value = 74995
$ws.Range("N123").Value = -84795
$ws.Range("H126").Value = 2437.4285
$ws.Range("I126").Value = 2877.6
$ws.Range("K126").Value = 8632.799999999999
$ws.Range("M126").Value = -8632.799999999999
